$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 10000
$ws.Range("I106").Value = 10000
$ws.Range("K106").Value = 10000
$ws.Range("M106").Value = -9369
$ws.Range("H129").Value = 3328.125
$ws.Range("I129").Value = 759
$ws.Range("J129").Value = 5897.25
$ws.Range("K129").Value = 2277
$ws.Range("L129").Value = 17691.75
$ws.Range("M129").Value = 2723
$ws.Range("N129").Value = -27691.75
$ws.Range("H132").Value = 1884.1875
$ws.Range("I132").Value = 1582.3572
$ws.Range("J132").Value = 3997
$ws.Range("K132").Value = 4747.071599999999
$ws.Range("L132").Value = 11991
$ws.Range("M132").Value = -2217.071599999999
$ws.Range("N132").Value = -17051
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
$ws.Range("H138").Value = 3779.3333
$ws.Range("I138").Value = 2848.5
$ws.Range("K138").Value = 8545.5
$ws.Range("M138").Value = -3405.5

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3531.0908
$ws.Range("I45").Value = 2495.5
$ws.Range("K45").Value = 2495.5
$ws.Range("M45").Value = -2118.5
$ws.Range("H122").Value = 1530.3334
$ws.Range("I122").Value = 1256.8572
$ws.Range("K122").Value = 3770.5716
$ws.Range("M122").Value = -1320.5716
$ws.Range("H132").Value = 1840.2667
$ws.Range("I132").Value = 1686
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 5058
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -2528
$ws.Range("N132").Value = -17060

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 3674.8215
$ws.Range("J7").Value = 81.44444
$ws.Range("L7").Value = 81.44444
$ws.Range("N7").Value = -307.44444
$ws.Range("H22").Value = 451.375
$ws.Range("I22").Value = 432
$ws.Range("J22").Value = 483.66666
$ws.Range("K22").Value = 432
$ws.Range("L22").Value = 483.66666
$ws.Range("M22").Value = -82
$ws.Range("N22").Value = -1183.66666
$ws.Range("H31").Value = 6580.9653
$ws.Range("I31").Value = 2754.125
$ws.Range("J31").Value = 8038.8096
$ws.Range("K31").Value = 2754.125
$ws.Range("L31").Value = 8038.8096
$ws.Range("M31").Value = -2459.125
$ws.Range("N31").Value = -8628.809600000001
$ws.Range("H34").Value = 6580.9653
$ws.Range("I34").Value = 2754.125
$ws.Range("J34").Value = 8038.8096
$ws.Range("K34").Value = 2754.125
$ws.Range("L34").Value = 8038.8096
$ws.Range("M34").Value = -2552.125
$ws.Range("N34").Value = -8442.809600000001
$ws.Range("H58").Value = 1864.5
$ws.Range("I58").Value = 835.8570999999999
$ws.Range("K58").Value = 835.8570999999999
$ws.Range("M58").Value = -632.8570999999999
$ws.Range("H132").Value = 2345.4443
$ws.Range("I132").Value = 2345.4443
$ws.Range("K132").Value = 7036.3329
$ws.Range("M132").Value = -4506.3329
$ws.Range("H136").Value = 1864.5
$ws.Range("I136").Value = 835.8570999999999
$ws.Range("K136").Value = 2507.5713
$ws.Range("M136").Value = 42.42870000000039

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 642.15
$ws.Range("I113").Value = 523.46155
$ws.Range("J113").Value = 862.5714
$ws.Range("K113").Value = 1570.38465
$ws.Range("L113").Value = 2587.7142
$ws.Range("M113").Value = 599.61535
$ws.Range("N113").Value = -6927.7142
$ws.Range("H140").Value = 2931.2856
$ws.Range("I140").Value = 2586.5
$ws.Range("K140").Value = 7759.5
$ws.Range("M140").Value = -2579.5
$ws.Range("H141").Value = 3045.8
$ws.Range("I141").Value = 1307.25
$ws.Range("K141").Value = 3921.75
$ws.Range("M141").Value = 1258.25

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("H80").Value = 7332
$ws.Range("H83").Value = 7332
$ws.Range("H102").Value = 1189.2
$ws.Range("I102").Value = 1189.2
$ws.Range("K102").Value = 1189.2
$ws.Range("M102").Value = 432.8
$ws.Range("H122").Value = 3288.8572
$ws.Range("I122").Value = 3637.2
$ws.Range("J122").Value = 2418
$ws.Range("K122").Value = 10911.6
$ws.Range("L122").Value = 7254
$ws.Range("M122").Value = -8461.599999999999
$ws.Range("N122").Value = -12154

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5098.5713
$ws.Range("J46").Value = 5616.6665
$ws.Range("L46").Value = 5616.6665
$ws.Range("N46").Value = -5992.6665
$ws.Range("H82").Value = 7499.8887
$ws.Range("I82").Value = 5750
$ws.Range("J82").Value = 7999.857
$ws.Range("K82").Value = 5750
$ws.Range("L82").Value = 7999.857
$ws.Range("M82").Value = -5389
$ws.Range("N82").Value = -8721.857
$ws.Range("H85").Value = 7499.8887
$ws.Range("I85").Value = 5750
$ws.Range("J85").Value = 7999.857
$ws.Range("K85").Value = 5750
$ws.Range("L85").Value = 7999.857
$ws.Range("M85").Value = -4502
$ws.Range("N85").Value = -10495.857
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H132").Value = 5148.8335
$ws.Range("I132").Value = 4978.6
$ws.Range("K132").Value = 14935.8
$ws.Range("M132").Value = -12405.8
$ws.Range("H136").Value = 6599.2
$ws.Range("I136").Value = 4000
$ws.Range("J136").Value = 7249
$ws.Range("K136").Value = 12000
$ws.Range("L136").Value = 21747
$ws.Range("M136").Value = -9450
$ws.Range("N136").Value = -26847

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4250.357
$ws.Range("I122").Value = 4292.1665
$ws.Range("K122").Value = 12876.4995
$ws.Range("M122").Value = -10426.4995
$ws.Range("H126").Value = 6739.533
$ws.Range("I126").Value = 4338.6
$ws.Range("J126").Value = 7940
$ws.Range("K126").Value = 13015.8
$ws.Range("L126").Value = 23820
$ws.Range("M126").Value = -10545.8
$ws.Range("N126").Value = -28760
$ws.Range("H131").Value = 70972.63
$ws.Range("J131").Value = 70972.63
$ws.Range("L131").Value = 70972.63
$ws.Range("N131").Value = -81052.63
$ws.Range("H136").Value = 3264.6667
$ws.Range("I136").Value = 2475.7778
$ws.Range("K136").Value = 7427.3334
$ws.Range("M136").Value = -4877.3334

